$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 1.95
$ws.Range("K2").Value = 2.3
$ws.Range("L2").Value = 5.6
$ws.Range("P2").Value = 3.96
$ws.Range("U2").Value = 1.74
$ws.Range("V2").Value = 2.04
$ws.Range("W2").Value = 6.6
$ws.Range("X2").Value = 6.6
$ws.Range("AA2").Value = 9.5
$ws.Range("AB2").Value = 18
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 50
$ws.Range("AH2").Value = 13.5
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 15
$ws.Range("AL2").Value = 45
$ws.Range("AM2").Value = 40
$ws.Range("AO2").Value = 6.8
$ws.Range("AP2").Value = 14.5
$ws.Range("AQ2").Value = 19.5
$ws.Range("AR2").Value = 40
$ws.Range("AT2").Value = 3.1
$ws.Range("AW2").Value = 7.5
$ws.Range("AX2").Value = 35
$ws.Range("AY2").Value = 35
$ws.Range("BB2").Value = 400
